# no-op test
$d = $word.ActiveDocument
Write-Output "paragraphs: $($d.Paragraphs.Count)"
